# Applies the cell-value changes for rows 8-25 and the Grand Total / Net
# Payable rows (27 & 29) of the "Bill Summary" sheet.
#
# Cells in columns A, D, E, G (and H/G on rows 27/29) are text-typed in the
# original workbook (t="str") even when their contents look like numbers
# (e.g. "4", "17", "56270.00"). Excel auto-converts a plain numeric-looking
# string assigned to Range.Value into a real Number, so those values are
# written with a leading single-quote (the standard Excel "treat as text"
# prefix) to keep them as text, matching the source file. A lone "'" is used
# for cells whose new value is an empty string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 38

# Row 9
$ws.Range("C9").Value = 85
$ws.Range("D9").Value = "'4"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
$ws.Range("G9").Value = "'56270.00"

# Row 10
$ws.Range("C10").Value = 47
$ws.Range("D10").Value = "'6"
$ws.Range("E10").Value = "On board"
$ws.Range("F10").Value = 136
$ws.Range("G10").Value = "'6392.00"

# Row 11
$ws.Range("A11").Value = "Each"
$ws.Range("C11").Value = 60
$ws.Range("D11").Value = "'4.0"
$ws.Range("E11").Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = "'3000.00"

# Row 12
$ws.Range("C12").Value = 35
$ws.Range("D12").Value = "'5.0"
$ws.Range("E12").Value = "Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F12").Value = 33
$ws.Range("G12").Value = "'1155.00"

# Row 13
$ws.Range("C13").Value = 70
$ws.Range("D13").Value = "'6.0"
$ws.Range("E13").Value = "Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F13").Value = 78
$ws.Range("G13").Value = "'5460.00"

# Row 14
$ws.Range("C14").Value = 100
$ws.Range("G14").Value = "'21900.00"

# Row 15
$ws.Range("A15").Value = "'"
$ws.Range("C15").Value = 76
$ws.Range("D15").Value = "'11.0"
$ws.Range("E15").Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = "'0.00"

# Row 16
$ws.Range("C16").Value = 44
$ws.Range("D16").Value = "'16"
$ws.Range("E16").Value = "20 mm"
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = "'1760.00"

# Row 17
$ws.Range("A17").Value = "R. mtr."
$ws.Range("C17").Value = 55
$ws.Range("D17").Value = "'17"
$ws.Range("E17").Value = "25 mm"
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = "'3080.00"

# Row 18
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = "'14.0"
$ws.Range("E18").Value = "Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "

# Row 19
$ws.Range("A19").Value = "Mtr."
$ws.Range("C19").Value = 52
$ws.Range("D19").Value = "'23"
$ws.Range("E19").Value = "8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = "'1040.00"

# Row 20
$ws.Range("C20").Value = 98
$ws.Range("D20").Value = "'15.0"
$ws.Range("E20").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."

# Row 21
$ws.Range("A21").Value = "'"
$ws.Range("C21").Value = 70
$ws.Range("D21").Value = "'29"
$ws.Range("E21").Value = "Single pole MCB   (With B/C curve tripping Characteristics)"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = "'0.00"

# Row 22
$ws.Range("A22").Value = "Each"
$ws.Range("C22").Value = 56
$ws.Range("D22").Value = "'32"
$ws.Range("E22").Value = " 50/63 A rating"
$ws.Range("F22").Value = 900
$ws.Range("G22").Value = "'50400.00"

# Row 23
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = "'18.0"
$ws.Range("E23").Value = "Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"

# Row 24
$ws.Range("C24").Value = 67
$ws.Range("D24").Value = "'34"
$ws.Range("E24").Value = "Metal door (single phase) IK-09 and IP-43 with Metal end box"

# Row 25
$ws.Range("C25").Value = 31

# Row 27
$ws.Range("G27").Value = "'150457.00"
$ws.Range("H27").Value = "'150457.00"

# Row 29
$ws.Range("G29").Value = "'150457.00"
$ws.Range("H29").Value = "'150457.00"
